$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns (D, E) store plain text (e.g. "22.454.54",
# "  +0.27%  ") rather than numbers. Setting .Value on a range whose
# NumberFormat is still "General" lets Excel auto-coerce numeric-looking
# text into a real number, so force Text format first, then restore the
# default "Normal" style afterwards so no stray formatting is left behind.
$affected = $ws.Range("D2:E51")
$affected.NumberFormat = "@"

$ws.Range("D2").Value = '22.454.54'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '1.573.58'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("D6").Value = '291.77'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").Value = '0.3729'
$ws.Range("D8").Value = '49.97'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.3396'
$ws.Range("E9").Value = '  -0.77%  '
$ws.Range("D10").Value = '0.07571'
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("D11").Value = '1.146'
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '21.30'
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("D14").Value = '6.021'
$ws.Range("D15").Value = '6.966'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").Value = '1.574.55'
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("E17").Value = '  -0.75%  '
$ws.Range("D18").Value = '91.00'
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("D19").Value = '0.06755'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '6.308'
$ws.Range("D22").Value = '16.31'
$ws.Range("E22").Value = '  -2.82%  '
$ws.Range("D23").Value = '12.17'
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("D24").Value = '22.442.69'
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("D25").Value = '2.342'
$ws.Range("E25").Value = '  -2.28%  '
$ws.Range("D26").Value = '2.691'
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("D27").Value = '20.08'
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").Value = '148.63'
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("D29").Value = '5.026'
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '125.63'
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("D31").Value = '1.750.53'
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("D32").Value = '1.061'
$ws.Range("E32").Value = '  +8.10%  '
$ws.Range("D33").Value = '6.189'
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("D34").Value = '1.988'
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("D35").Value = '9.864'
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("D36").Value = '0.08382'
$ws.Range("E36").Value = '  -1.12%  '
$ws.Range("D37").Value = '0.02496'
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("D39").Value = '0.2304'
$ws.Range("E39").Value = '  -0.61%  '
$ws.Range("D40").Value = '0.06521'
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").Value = '5.476'
$ws.Range("E41").Value = '  +0.94%  '
$ws.Range("D42").Value = '11.32'
$ws.Range("E42").Value = '  -1.26%  '
$ws.Range("D43").Value = '0.6224'
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '13.97'
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").Value = '3.814'
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("D47").Value = '0.5807'
$ws.Range("E47").Value = '  -2.79%  '
$ws.Range("D48").Value = '129.75'
$ws.Range("E48").Value = '  +3.40%  '
$ws.Range("D49").Value = '2.069'
$ws.Range("E49").Value = '  -0.99%  '
$ws.Range("D50").Value = '1.223'
$ws.Range("E50").Value = '  -5.57%  '
$ws.Range("E51").Value = '  -0.01%  '

$affected.Style = "Normal"
